$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so that numeric-looking
# strings (e.g. "0.603", "64.274.76") are preserved exactly as text,
# matching the original inlineStr cell contents instead of being
# auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.274.76'
$ws.Range("E2").Value = '  -2.85%  '

$ws.Range("D3").Value = '3.160.64'
$ws.Range("E3").Value = '  -4.62%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '570.29'
$ws.Range("E5").Value = '  -2.56%  '

$ws.Range("D6").Value = '168.03'
$ws.Range("E6").Value = '  -6.87%  '

$ws.Range("D7").Value = '0.603'
$ws.Range("E7").Value = '  -7.42%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '3.184.95'
$ws.Range("E9").Value = '  -3.83%  '

$ws.Range("E10").Value = '  -5.69%  '

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -3.22%  '

$ws.Range("D13").Value = '3.713.73'
$ws.Range("E13").Value = '  -4.61%  '

$ws.Range("E14").Value = '  -1.46%  '

$ws.Range("D15").Value = '64.408.34'
$ws.Range("E15").Value = '  -2.71%  '

$ws.Range("D16").Value = '25.32'
$ws.Range("E16").Value = '  -4.57%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0000157'
$ws.Range("E17").Value = '  -4.14%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.157.88'
$ws.Range("E18").Value = '  -6.33%  '

$ws.Range("D19").Value = '417.25'
$ws.Range("E19").Value = '  -1.55%  '

$ws.Range("D20").Value = '12.84'
$ws.Range("E20").Value = '  -1.82%  '

$ws.Range("D21").Value = '5.30'
$ws.Range("E21").Value = '  -3.65%  '

$ws.Range("D22").Value = '7.14'
$ws.Range("E22").Value = '  -3.09%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").Value = '69.55'
$ws.Range("E25").Value = '  -2.96%  '

$ws.Range("D26").Value = '0.203'
$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("D27").Value = '0.498'
$ws.Range("E27").Value = '  -2.93%  '

$ws.Range("D28").Value = '0.0000103'
$ws.Range("E28").Value = '  -10.02%  '

$ws.Range("E29").Value = '  -2.83%  '

$ws.Range("E30").Value = '  -0.13%  '

$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  -4.98%  '

$ws.Range("D32").Value = '21.70'
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").Value = '5.05'
$ws.Range("E34").Value = '  -2.23%  '

$ws.Range("D35").Value = '6.35'
$ws.Range("E35").Value = '  -3.98%  '

$ws.Range("D36").Value = '1.12'
$ws.Range("E36").Value = '  -5.15%  '

$ws.Range("D37").Value = '155.13'
$ws.Range("E37").Value = '  -2.99%  '

$ws.Range("E38").Value = '  -5.16%  '

$ws.Range("D39").Value = '2.704.06'
$ws.Range("E39").Value = '  -5.43%  '

$ws.Range("E40").Value = '  -5.70%  '

$ws.Range("D41").Value = '4.21'

$ws.Range("D42").Value = '23.90'
$ws.Range("E42").Value = '  -9.19%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '0.717'
$ws.Range("E43").Value = '  -5.47%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '39.00'
$ws.Range("E44").Value = '  -1.72%  '

$ws.Range("D45").Value = '0.0618'
$ws.Range("E45").Value = '  -6.32%  '

$ws.Range("D46").Value = '5.49'
$ws.Range("E46").Value = '  -7.17%  '

$ws.Range("E47").Value = '  -3.31%  '

$ws.Range("D48").Value = '21.34'
$ws.Range("E48").Value = '  -7.40%  '

$ws.Range("D49").Value = '288.51'
$ws.Range("E49").Value = '  -6.95%  '

$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").Value = '0.0990'
$ws.Range("E51").Value = '  -5.27%  '

# Restore the default (Normal) style so the cells keep no explicit
# style index, matching the original workbook formatting.
$dataRange.Style = "Normal"
